$wb = $excel.ActiveWorkbook

$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "CUMPLIMIENTO MENSUAL"

$ws.Range("A1").Value = "ASESOR"
$ws.Range("B1").Value = "GRUPO"
$ws.Range("C1").Value = "PRESUPUESTO"
$ws.Range("D1").Value = "VENTA"
$ws.Range("E1").Value = "POR CUMPLIR"
$ws.Range("F1").Value = "CUMPLIMIENTO"

$asesor = "ALMEIDA CUATIN JHONATHANN CARLOS"

$groups = @(
  @("240X120 PORCELANATO", 260.285000070615),
  @("240X80 PORCELANATO", 3120.1145),
  @("FREGADEROS DE COCINA", 646.361575487259),
  @("GRANITO", 238.32),
  @("GRIFERIAS", 106.82),
  @("INODOROS", 1260),
  @("LAVABOS", 625),
  @("LED", 300),
  @("NO RESURTIBLES", 650.25),
  @("OTROS", 0),
  @("PANELES DECORATIVOS", 350),
  @("PANELES PU", 230),
  @("PANELES PVC", 483),
  @("PIEDRA SINTERIZADA", 527.03),
  @("PORCELANATO", 18798.61),
  @("PUERTAS DE SEGURIDAD", 342),
  @("SAL SOLUBLE", 1600)
)

$row = 2
foreach ($g in $groups) {
  $ws.Range("A$row").Value = $asesor
  $ws.Range("B$row").Value = $g[0]
  $ws.Range("C$row").Value = $g[1]
  $ws.Range("C$row").NumberFormat = "`"$`"#,##0.00"
  $ws.Range("D$row").Value = 0
  $ws.Range("D$row").NumberFormat = "`"$`"#,##0.00"
  $ws.Range("E$row").Value = $g[1]
  $ws.Range("E$row").NumberFormat = "`"$`"#,##0.00"
  $ws.Range("F$row").Value = 0
  $ws.Range("F$row").NumberFormat = "0.00%"
  $row = $row + 1
}

$totalRow = $row
$ws.Range("B$totalRow").Value = "TOTAL"
$ws.Range("B$totalRow").HorizontalAlignment = -4152
$ws.Range("C$totalRow").Value = 29537.79107555788
$ws.Range("C$totalRow").NumberFormat = "`"$`"#,##0.00"
$ws.Range("D$totalRow").Value = 0
$ws.Range("D$totalRow").NumberFormat = "`"$`"#,##0.00"
$ws.Range("E$totalRow").Value = 29537.79107555788
$ws.Range("E$totalRow").NumberFormat = "`"$`"#,##0.00"
$ws.Range("F$totalRow").Value = 0
$ws.Range("F$totalRow").NumberFormat = "0.00%"

$ws.Columns.Item(1).ColumnWidth = 34
$ws.Columns.Item(2).ColumnWidth = 22
$ws.Columns.Item(3).ColumnWidth = 22
$ws.Columns.Item(4).ColumnWidth = 11
$ws.Columns.Item(5).ColumnWidth = 22
$ws.Columns.Item(6).ColumnWidth = 18

Write-Output "done"
